# Generate Report for Handoff
# - Set Priority ("ht") for the affected rows on the zh-cn and de-de sheets
# - Bump the related timestamps (Overview "Latest HO Xliff Generate Date" and
#   the zh-cn / de-de "Latest Handoff Datetime") to reflect the new report run

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 11, 12, 13, 14)

# --- Overview sheet: Latest HO Xliff Generate Date (column G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-31 16:23:58"
}

# --- zh-cn sheet: Priority (column E) + Latest Handoff Datetime (column H) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-31 16:23:54"
}

# --- de-de sheet: Priority (column E) + Latest Handoff Datetime (column H) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 8).Value = "2016-08-31 16:23:58"
}
